$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 261
    3  = 262
    4  = 264
    5  = 267
    6  = 270
    7  = 272
    8  = 273
    9  = 275
    10 = 277
    11 = 279
    12 = 282
    13 = 284
    14 = 285
    15 = 46
    16 = 72
    17 = 78
    18 = 144
    19 = 161
    20 = 200
    21 = 231
    22 = 302
    23 = 318
    24 = 347
    25 = 427
    26 = 450
    27 = 497
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
